$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ D = "D2";  Dv = "28.121.31";    E = "E2";  Ev = "  -1.01%  " },
    @{ D = "D3";  Dv = "1.794.79";     E = "E3";  Ev = "  +0.00%  " },
    @{              E = "E4";  Ev = "  -0.12%  " },
    @{ D = "D5";  Dv = "317.52";       E = "E5";  Ev = "  +0.95%  " },
    @{ D = "D6";  Dv = "1.001";        E = "E6";  Ev = "  -0.11%  " },
    @{ D = "D7";  Dv = "0.5404";       E = "E7";  Ev = "  -0.42%  " },
    @{ D = "D8";  Dv = "0.3784" },
    @{ D = "D9";  Dv = "0.07448";      E = "E9";  Ev = "  -1.71%  " },
    @{ D = "D10"; Dv = "41.74";        E = "E10"; Ev = "  -1.70%  " },
    @{ D = "D11"; Dv = "1.093";        E = "E11"; Ev = "  -2.59%  " },
    @{              E = "E12"; Ev = "  -0.12%  " },
    @{ D = "D13"; Dv = "20.53";        E = "E13"; Ev = "  -2.70%  " },
    @{ D = "D14"; Dv = "6.114";        E = "E14"; Ev = "  -1.19%  " },
    @{ D = "D15"; Dv = "7.242";        E = "E15"; Ev = "  -2.18%  " },
    @{ D = "D16"; Dv = "1.786.83";     E = "E16"; Ev = "  -0.65%  " },
    @{ D = "D17"; Dv = "89.08";        E = "E17"; Ev = "  -2.90%  " },
    @{ D = "D18"; Dv = "0.00001060";   E = "E18"; Ev = "  -0.85%  " },
    @{ D = "D19"; Dv = "0.06484";      E = "E19"; Ev = "  +0.37%  " },
    @{              E = "E20"; Ev = "  -0.13%  " },
    @{ D = "D21"; Dv = "17.26";        E = "E21"; Ev = "  -0.44%  " },
    @{ D = "D22"; Dv = "5.905";        E = "E22"; Ev = "  -0.97%  " },
    @{ D = "D23"; Dv = "28.135.09";    E = "E23"; Ev = "  -0.99%  " },
    @{ D = "D24"; Dv = "11.16";        E = "E24"; Ev = "  -1.67%  " },
    @{ D = "D25"; Dv = "2.088";        E = "E25"; Ev = "  -1.52%  " },
    @{ D = "D26"; Dv = "155.08";       E = "E26"; Ev = "  -2.96%  " },
    @{ D = "D27"; Dv = "20.27";        E = "E27"; Ev = "  -2.00%  " },
    @{ D = "D28"; Dv = "1.998.20";     E = "E28"; Ev = "  -0.35%  " },
    @{ D = "D29"; Dv = "2.283";        E = "E29"; Ev = "  -4.71%  " },
    @{ D = "D30"; Dv = "121.12";       E = "E30"; Ev = "  -1.72%  " },
    @{ D = "D31"; Dv = "1.119";        E = "E31"; Ev = "  +0.00%  " },
    @{ D = "D32"; Dv = "0.1057";       E = "E32"; Ev = "  +3.50%  " },
    @{ D = "D33"; Dv = "3.654";        E = "E33"; Ev = "  -1.18%  " },
    @{ D = "D34"; Dv = "5.553";        E = "E34"; Ev = "  -3.19%  " },
    @{ D = "D35"; Dv = "0.2258";       E = "E35"; Ev = "  -2.95%  " },
    @{ D = "D36"; Dv = "0.06492";      E = "E36"; Ev = "  +1.86%  " },
    @{ D = "D37"; Dv = "0.02292";      E = "E37"; Ev = "  -1.16%  " },
    @{              E = "E38"; Ev = "  -2.39%  " },
    @{ D = "D39"; Dv = "8.453";        E = "E39"; Ev = "  -3.82%  " },
    @{              E = "E40"; Ev = "  +4.42%  " },
    @{ D = "D41"; Dv = "0.6179";       E = "E41"; Ev = "  -3.39%  " },
    @{ D = "D42"; Dv = "11.08";        E = "E42"; Ev = "  -4.58%  " },
    @{ D = "D43"; Dv = "1.173" },
    @{ D = "D44"; Dv = "0.9999";       E = "E44"; Ev = "  -0.12%  " },
    @{ D = "D45"; Dv = "13.28";        E = "E45"; Ev = "  -2.35%  " },
    @{ D = "D46"; Dv = "3.678";        E = "E46"; Ev = "  +0.09%  " },
    @{ D = "D47"; Dv = "0.5781";       E = "E47"; Ev = "  -3.11%  " },
    @{ D = "D48"; Dv = "124.21";       E = "E48"; Ev = "  -1.52%  " },
    @{ D = "D49"; Dv = "1.191";        E = "E49"; Ev = "  +3.58%  " },
    @{ D = "D50"; Dv = "1.923";        E = "E50"; Ev = "  -3.06%  " },
    @{ D = "D51"; Dv = "0.06814";      E = "E51"; Ev = "  -1.11%  " }
)

# Force the Price column cells being written to stay as plain text, matching
# the workbook's existing convention (inline/shared strings), instead of
# letting Excel auto-convert numeric-looking text into real numbers.
foreach ($row in $updates) {
    if ($row.ContainsKey("D")) {
        $ws.Range($row.D).NumberFormat = "@"
    }
}

foreach ($row in $updates) {
    if ($row.ContainsKey("D")) {
        $ws.Range($row.D).Value = $row.Dv
    }
    if ($row.ContainsKey("E")) {
        $ws.Range($row.E).Value = $row.Ev
    }
}

# Restore the default "Normal" style on the touched Price cells so that no
# stray number-format style lingers on them (matches original: no "s" attr).
foreach ($row in $updates) {
    if ($row.ContainsKey("D")) {
        $ws.Range($row.D).Style = "Normal"
    }
}
